$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ver-Construção1")

$ws.Range("D16").Value = "Sim"
$ws.Range("D18").Value = "Sim"
$ws.Range("D19").Value = "Sim"
$ws.Range("D21").Value = "NA"
$ws.Range("D23").Value = "NA"
$ws.Range("D24").Value = "Sim"
$ws.Range("D25").Value = "Sim"
$ws.Range("D26").Value = "Sim"
$ws.Range("D27").Value = "NA"
$ws.Range("D28").Value = "NA"
$ws.Range("D30").Value = "Sim"
$ws.Range("D31").Value = "Sim"
$ws.Range("D32").Value = "Sim"
$ws.Range("D33").Value = "Sim"
$ws.Range("D34").Value = "NA"
$ws.Range("D35").Value = "NA"
$ws.Range("D36").Value = "NA"
$ws.Range("D38").Value = "Sim"
$ws.Range("D39").Value = "Sim"
$ws.Range("D40").Value = "NA"
$ws.Range("D41").Value = "Sim"
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = "Sim"
$ws.Range("D45").Value = "Sim"
$ws.Range("D46").Value = "Sim"
$ws.Range("D47").Value = "Sim"
$ws.Range("D49").Value = "NA"
$ws.Range("D50").Value = "NA"
$ws.Range("D51").Value = "NA"

$wb.Save()
